# Update BOM unit-cost (G) and extended-cost (H) values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G5").Value = 0.62
$ws.Range("H5").Value = 1.24

$ws.Range("G6").Value = 0.03271
$ws.Range("H6").Value = 0.32705

$ws.Range("G7").Value = 0.0257
$ws.Range("H7").Value = 0.257

$ws.Range("G11").Value = 0.094
$ws.Range("H11").Value = 0.94

$ws.Range("G14").Value = 0.0687
$ws.Range("H14").Value = 0.0687

$ws.Range("G18").Value = 0.3358
$ws.Range("H18").Value = 0.3358

$ws.Range("G26").Value = 0.0271
$ws.Range("H26").Value = 0.0542

$ws.Range("G30").Value = 0.0369
$ws.Range("H30").Value = 0.2214

$ws.Range("G38").Value = 0.098
$ws.Range("H38").Value = 0.098
